$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, reusing the same header style as the
# other header cells (e.g. G1) by copying its formatting.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the Save values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
